# Medidas.xlsx edit: add VPN and f1-score Negativo rows, "+"/"-" indicator
# column, and a new "NPB" section label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the two new rows (this shifts the rows below them down by one
#    each time, exactly like using the right click "Insert" on a row
#    header in Excel).
# ---------------------------------------------------------------------
$ws.Rows("6:6").Insert()
$ws.Rows("8:8").Insert()

# ---------------------------------------------------------------------
# 2) Fill in the new VPN row (row 6)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "VPN"
$ws.Range("C6").Value = "0.8850514105"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "mide la proporción de verdaderos negativos sobre el total de predicciones negativas realizadas por el modelo"

# ---------------------------------------------------------------------
# 3) Fill in the new f1-score Negativo row (row 8)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "f1-score Negativo"
$ws.Range("C8").Value = "0.8936272160996646"
$ws.Range("E8").Value = "-"

# ---------------------------------------------------------------------
# 4) Add the "+"/"-" indicator column for the remaining metric rows
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "+"
$ws.Range("E4").Value = "-"
$ws.Range("E5").Value = "+"
$ws.Range("E7").Value = "+"

# ---------------------------------------------------------------------
# 5) Style fix-ups so the new cells match the look of the rest of the
#    table.
# ---------------------------------------------------------------------

# A6 / A8 should look like the other metric-name cells in column A
# (already carried over from the row insert, but make sure explicitly).
$ws.Range("A6,A8").Font.Name = "Calibri"
$ws.Range("A6,A8").Font.Size = 11
$ws.Range("A6,A8").Interior.ThemeColor = 8
$ws.Range("A6,A8").Interior.TintAndShade = 0.6

# C6 keeps the value-cell look (font + fill) but without a border.
$ws.Range("C6").Font.Name = "Var(--jp-code-font-family)"
$ws.Range("C6").Font.Size = 10
$ws.Range("C6").HorizontalAlignment = -4131
$ws.Range("C6").VerticalAlignment = -4108
$ws.Range("C6").Interior.ThemeColor = 9
$ws.Range("C6").Interior.TintAndShade = 0.6
$ws.Range("C6").Borders.LineStyle = 0

# C8 keeps the normal bordered value-cell look.
$ws.Range("C8").Font.Name = "Var(--jp-code-font-family)"
$ws.Range("C8").Font.Size = 10
$ws.Range("C8").HorizontalAlignment = -4131
$ws.Range("C8").VerticalAlignment = -4108
$ws.Range("C8").Interior.ThemeColor = 9
$ws.Range("C8").Interior.TintAndShade = 0.6
$ws.Range("C8").Borders.LineStyle = 1

# D6/D8 stay empty, formatted like the rest of column D.
$ws.Range("D6,D8").Font.Name = "Var(--jp-code-font-family)"
$ws.Range("D6,D8").Font.Size = 10
$ws.Range("D6,D8").HorizontalAlignment = -4131
$ws.Range("D6,D8").VerticalAlignment = -4108
$ws.Range("D6,D8").Borders.LineStyle = 1

# F6 gets the regular Segoe UI description-cell look.
$ws.Range("F6").Font.Name = "Segoe UI"
$ws.Range("F6").Font.Size = 11

# F8 is left empty, but marked with an underlined font (placeholder for a
# description still to be written).
$ws.Range("F8").Font.Name = "Segoe UI"
$ws.Range("F8").Font.Size = 11
$ws.Range("F8").Font.Underline = 2

# ---------------------------------------------------------------------
# 6) New "NPB" section label two rows below the table (row 16, leaving
#    row 15 blank).
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "NPB"
$ws.Range("A16").Font.Name = "Calibri"
$ws.Range("A16").Font.Size = 11
$ws.Range("A16").Interior.ThemeColor = 8
$ws.Range("A16").Interior.TintAndShade = 0.6
$ws.Range("A16").Borders.LineStyle = 0

# ---------------------------------------------------------------------
# 7) Match the final selection shown in the workbook.
# ---------------------------------------------------------------------
$ws.Range("F8").Select()

Write-Output "done"
